$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39 was the last row and used the "last row" date style (s=3, YYYY-MM-DD).
# Now that a new last row is appended, row 39 reverts to the regular
# date-time style (s=2) and the new row 40 takes on the "last row" style (s=3).
$ws.Range("A39").NumberFormat = $ws.Range("A2").NumberFormat

$ws.Range("A40").Value = 45625
$ws.Range("B40").Value = 106
$ws.Range("C40").Value = 87
$ws.Range("D40").Value = 95

$ws.Range("A40").NumberFormat = "YYYY-MM-DD"
